# Re-parse data & regenerate all survey-answer sheets
# Rewrites the q10_lang_subtasks answer rows (sorted by donor_id) with a
# refreshed timestamp, growing the sheet from 8 data rows to 23 data rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$timestamp = 45854.65518179266
$question  = "q10_lang_subtasks"
$noChoice  = "I did not choose “Language practice or translation”"
$allFive   = "Translating full texts between languages|Improving grammar or style in a target language|Vocabulary drills or word lists|Conversational practice – dialogue role-play|Pronunciation or phonetic guidance"
$fourNoTranslate = "Improving grammar or style in a target language|Vocabulary drills or word lists|Conversational practice – dialogue role-play|Pronunciation or phonetic guidance"
$fourNoPron      = "Translating full texts between languages|Improving grammar or style in a target language|Vocabulary drills or word lists|Conversational practice – dialogue role-play"
$threeMid        = "Improving grammar or style in a target language|Vocabulary drills or word lists|Conversational practice – dialogue role-play"

$rows = @(
    @("0ce5dd49", $noChoice),
    @("2c1001cb", $noChoice),
    @("37cc37bf", $fourNoTranslate),
    @("43faa0b9", $allFive),
    @("4abe3e88", $fourNoPron),
    @("50164f59", $noChoice),
    @("5cf70f79", $allFive),
    @("5da96769", $noChoice),
    @("6ca3e2f6", $noChoice),
    @("790a4fcb", $threeMid),
    @("802cc63a", $noChoice),
    @("85c3ea4d", $allFive),
    @("942dfafb", $noChoice),
    @("9bc6ba8c", $noChoice),
    @("a2d65af2", $noChoice),
    @("a46f1771", $allFive),
    @("ad58f9da", $allFive),
    @("c7d9a301", $allFive),
    @("ce8732ff", $noChoice),
    @("d6f1d567", $noChoice),
    @("da9326c9", $fourNoPron),
    @("e09ca7bf", $noChoice),
    @("ef53a641", $noChoice)
)

for ($i = 0; $i -lt $rows.Length; $i++) {
    $r = $i + 2
    $donorId = $rows[$i][0]
    $category = $rows[$i][1]

    $ws.Cells.Item($r, 1).Value2 = $donorId
    $ws.Cells.Item($r, 2).Value2 = $category
    $ws.Cells.Item($r, 3).Value2 = $question
    $ws.Cells.Item($r, 4).Value2 = $timestamp
    $ws.Cells.Item($r, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
}

Write-Host "Wrote $($rows.Length) data rows (A2:D$($rows.Length + 1))"
